$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = $origStyle
}

Set-TextValue 4 3 "356"
Set-TextValue 4 4 "1296368.92"
Set-TextValue 9 3 "54"
Set-TextValue 9 4 "129428.41"
Set-TextValue 10 3 "329"
Set-TextValue 10 4 "1069996.74"
Set-TextValue 11 3 "135"
Set-TextValue 11 4 "557891.77"
Set-TextValue 12 3 "32"
Set-TextValue 12 4 "152120.00"
Set-TextValue 13 3 "6"
Set-TextValue 13 4 "32000.00"
Set-TextValue 14 3 "15"
Set-TextValue 14 4 "38000.00"
Set-TextValue 15 3 "98"
Set-TextValue 15 4 "254152.38"
Set-TextValue 16 3 "413"
Set-TextValue 16 4 "1262472.14"
Set-TextValue 18 3 "45"
Set-TextValue 18 4 "226045.00"
Set-TextValue 33 3 "100"
Set-TextValue 33 4 "268748.00"
Set-TextValue 34 3 "524"
Set-TextValue 34 4 "1633546.26"
Set-TextValue 37 3 "26"
Set-TextValue 37 4 "170500.00"
Set-TextValue 50 3 "95"
Set-TextValue 50 4 "269768.17"
Set-TextValue 51 3 "556"
Set-TextValue 51 4 "1869799.52"
Set-TextValue 52 3 "255"
Set-TextValue 52 4 "1086446.76"
Set-TextValue 53 3 "85"
Set-TextValue 53 4 "489878.23"
Set-TextValue 54 3 "24"
Set-TextValue 54 4 "133213.00"
Set-TextValue 55 3 "17"
Set-TextValue 55 4 "54720.65"
Set-TextValue 97 3 "276"
Set-TextValue 97 4 "714279.43"
Set-TextValue 98 3 "1178"
Set-TextValue 98 4 "3534360.70"
Set-TextValue 99 3 "440"
Set-TextValue 99 4 "1742594.02"
Set-TextValue 100 3 "116"
Set-TextValue 100 4 "527000.00"
Set-TextValue 102 3 "65"
Set-TextValue 102 4 "141000.00"
